$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the browser name used in the step
$ws.Range("F2").Value = "FIREFOX"

# Move the active selection (matches the recorded cursor position in the diff)
$ws.Range("F9").Select()
